$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C = y_0_forecast, Column E = y_1_forecast
# Update rows 2-19 with new forecast values (bugfixed evaluation / simulated rt_data)

$values = @(
    @{Row=2;  C=6.056254825277896;   E=4.566338461218011},
    @{Row=3;  C=8.604123301398037;   E=8.260999835306727},
    @{Row=4;  C=5.436647924209592;   E=4.862860110364875},
    @{Row=5;  C=6.334380382529425;   E=5.829578861489648},
    @{Row=6;  C=3.88993859232436;    E=3.4300351921007},
    @{Row=7;  C=2.513767348245044;   E=2.479992751939486},
    @{Row=8;  C=2.723916849952834;   E=2.551173534479334},
    @{Row=9;  C=1.326505206336948;   E=1.211929054838756},
    @{Row=10; C=2.321003614014883;   E=2.610201636760778},
    @{Row=11; C=2.468891199411116;   E=2.593292206016984},
    @{Row=12; C=3.133596157287766;   E=3.419422858788335},
    @{Row=13; C=0.08486825492834971; E=0.9311475558545057},
    @{Row=14; C=2.405224065057476;   E=2.152035263856344},
    @{Row=15; C=-0.678826357714013;  E=0.02414656897629097},
    @{Row=16; C=0.7010162698181555;  E=-0.6730402944081559},
    @{Row=17; C=2.158153176293576;   E=1.298013848993262},
    @{Row=18; C=-0.1645072558042915; E=0.1619933518385297},
    @{Row=19; C=3.131832690451031;   E=2.199676451050503}
)

foreach ($item in $values) {
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}
